$d = $word.ActiveDocument
$d.Content.Find.Execute("$20.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "$10.00", 2)
